$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.822.89"
$ws.Range("E2").Value = "  -0.66%  "

# Row 3
$ws.Range("D3").Value = "2.499.94"
$ws.Range("E3").Value = "  +2.86%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.33%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.15%  "

# Row 9
$ws.Range("D9").Value = "2.524.09"
$ws.Range("E9").Value = "  +3.53%  "

# Row 10
$ws.Range("E10").Value = "  +1.30%  "

# Row 11
$ws.Range("E11").Value = "  +0.37%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.55%  "

# Row 14
$ws.Range("D14").Value = "2.939.23"
$ws.Range("E14").Value = "  +2.85%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.22%  "

# Row 16
$ws.Range("D16").Value = "58.756.60"
$ws.Range("E16").Value = "  -0.64%  "

# Row 17
$ws.Range("E17").Value = "  +1.45%  "

# Row 18
$ws.Range("D18").Value = "2.517.27"
$ws.Range("E18").Value = "  +1.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.30%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.46%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.63%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.437"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.45%  "

# Row 26
$ws.Range("E26").Value = "  +1.64%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.44%  "

# Row 28
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.613.61"
$ws.Range("E28").Value = "  +3.00%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0770"
$ws.Range("E30").Value = "  +0.58%  "

# Row 31
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "

# Row 33
$ws.Range("E33").Value = "  -6.63%  "

# Row 34
$ws.Range("E34").Value = "  -0.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.77%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.71%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.28%  "

# Row 39
$ws.Range("E39").Value = "  -7.97%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.87%  "

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "298.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.69%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "

# Row 44
$ws.Range("B44").Value = "SuiNetwork"
$ws.Range("C44").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.822"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.600"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "  +0.31%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0928"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0227"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
